# Auto-generated edit script implementing the QA Round 2 diff.
$wb = $excel.ActiveWorkbook

# --- Step 1: split 'cumcontrol' into 'cumcontrol1' + 'cumcontrol2' ---
$cc = $wb.Worksheets.Item("cumcontrol")
# Duplicate the sheet first (keeps identical formatting/column widths/styles);
# the copy is placed immediately after the original.
$cc.Copy($null, $cc)
$cc2 = $wb.Worksheets.Item("cumcontrol (2)")
$cc2.Name = "cumcontrol2"
$cc.Name = "cumcontrol1"

# --- cumcontrol1: update text/notes, row names (A) stay the same ---
$cc.Range("B2").Value = 'just hold on a little more, I want the last thing you see to be this 🥺'
$cc.Range("B3").Value = 'wait for me... I have one more thing and I want you to see it before we finish'
$cc.Range("C3").Value = 'DELAY. Send PPV.'
$cc.Range("B4").Value = 'stay with me, I''m almost there too... watch this'
$cc.Range("C4").Value = 'SYNC variant. Send PPV.'
$cc.Range("B5").Value = 'I want us to finish together... open this and let go with me 🌸'
$cc.Range("C5").Value = 'SYNC. Send PPV.'
$cc.Range("B6").Value = 'please don''t finish yet... I''m not ready for this to be over 🥺'
$cc.Range("B7").Value = 'not yet babe... I want this to last a little longer with you 🌸'
$cc.Range("C7").Value = 'CONTROL.'

# --- cumcontrol2: update row names (A), text (B) and notes (C) ---
$cc2.Range("A2").Value = 'delay2'
$cc2.Range("B2").Value = 'just a little longer for me babe? the next one is special 🥺'
$cc2.Range("C2").Value = 'DELAY variant.'
$cc2.Range("A3").Value = 'delay1'
$cc2.Range("B3").Value = 'please wait... what I''m about to send, I want you to really take it in'
$cc2.Range("C3").Value = 'DELAY. Send PPV.'
$cc2.Range("A4").Value = 'sync2'
$cc2.Range("B4").Value = 'I need you to see this before we both let go'
$cc2.Range("C4").Value = 'SYNC variant.'
$cc2.Range("A5").Value = 'sync1'
$cc2.Range("B5").Value = 'okay babe... together, right now... open this 🌸'
$cc2.Range("C5").Value = 'SYNC. Send PPV.'
$cc2.Range("A6").Value = 'edge2'
$cc2.Range("B6").Value = 'don''t rush... this is too good to end yet'
$cc2.Range("C6").Value = 'EDGE variant.'
$cc2.Range("A7").Value = 'edge1'
$cc2.Range("B7").Value = 'slow down babe... I want to feel every second of this with you 💕'
$cc2.Range("C7").Value = 'CONTROL.'

# --- Step 2: remaining single-cell text edits across other sheets ---
$ws = $wb.Worksheets.Item('LanaJourney')
$ws.Range('B3').Value = 'wow that was... incredible 🌸'
$ws.Range('B8').Value = 'fuuuck 🥺'
$ws.Range('B9').Value = 'this is only for you, nobody else gets to see this 💕'
$ws.Range('B13').Value = 'I want you so bad right now 🌸'
$ws.Range('B14').Value = 'fuuuck 😳'
$ws.Range('B15').Value = 'okay I really want to show you something... look 💕'
$ws.Range('B19').Value = 'omg cutie you''re making me feel so many things right now... 🌸'
$ws.Range('B20').Value = 'tell me what you think 💕'
$ws.Range('B21').Value = 'I''ve never shown this to anyone before... this is just for you 🌸'
$ws.Range('B23').Value = 'you''re making me feel some type of way right now 🌸'
$ws.Range('B29').Value = 'like I keep getting these butterflies and I can''t explain it, you make me feel so comfortable, that''s kinda rare for me 💕'
$ws = $wb.Worksheets.Item('nosex1')
$ws.Range('B3').Value = 'I can''t help it, there''s something about you that''s making me feel things nobody else does 🥺'
$ws = $wb.Worksheets.Item('PersonalLana')
$ws.Range('B2').Value = 'I''m 18, just started college and I''m loving this whole new chapter of my life 🌸'
$ws = $wb.Worksheets.Item('PositiveSpin')
$ws.Range('B2').Value = 'honestly I love talking to someone who actually knows what they want, you make me feel so comfortable and that means everything to me 🌸'
$ws = $wb.Worksheets.Item('ReEngagement')
$ws.Range('B2').Value = 'so after we stopped talking I did something special and you''re the only one I want to share it with 💕'
$ws = $wb.Worksheets.Item('discount2')
$ws.Range('B5').Value = 'a discount? babe this is really personal to me and I don''t usually do this 🥺'
$ws = $wb.Worksheets.Item('boosters')
$ws.Range('B5').Value = 'I''ve literally never felt like this 🥺'

# --- sanity check: print final sheet order ---
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }

